# Re-test on Subsumes7: update the raw measured timing values (column P,
# rows 3-12) with the results of a fresh run, then let all the dependent
# formulas (Q/R/S columns, Samenvatting_6, Samenvatting_7, Sheet1, charts)
# recalculate automatically.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Subsumes7")

$newValues = @(24823921974, 29165320450, 24624222015, 28927545003, 24258287724, 29171312688, 29036512816, 29133569029, 29088751267, 28979078576)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 16).Value = $newValues[$i]
}

$excel.Calculate()

# Restore the selection on Subsumes7 to where it ended up after the re-test.
$ws.Range("Q14").Select()

# The active tab moved off the Grafiek6 chart sheet onto Samenvatting_7.
$wb.Worksheets.Item("Samenvatting_7").Activate()
$wb.Worksheets.Item("Samenvatting_7").Range("C9").Select()
